# Update the "Contenu du stage" statistics (counts + percentage labels)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Update the counts (column E) for the "Contenu du stage" breakdown
$ws.Range("E16").Value = 8
$ws.Range("E17").Value = 18
$ws.Range("E19").Value = 2

# Update the percentage labels (column G) to match the new counts.
# Force the cells to text format so the "%" strings are not
# reinterpreted as numeric percentages by Excel.
$ws.Range("G16:G19").NumberFormat = "@"
$ws.Range("G16").Value = "28.57 %"
$ws.Range("G17").Value = "64.29 %"
$ws.Range("G18").Value = "0 %"
$ws.Range("G19").Value = "7.14 %"

$wb.Save()
